# Generate Report for Handback
# Applies the "handback" status update to the localization-status workbook:
#  - Overview sheet: status columns (zh-cn / de-de) move from "In Translation"
#    to "Handed back: in sync with en-US"
#  - zh-cn / de-de detail sheets: fill in the "Latest Target File" and
#    "Latest Handback File" columns (with hyperlinks on the target file) and
#    stamp the "Latest Handback DateTime" column
#  - widen the columns that now hold the longer values

$wb = $excel.ActiveWorkbook

$ovr   = $wb.Worksheets.Item("Overview")
$zhcn  = $wb.Worksheets.Item("zh-cn")
$dede  = $wb.Worksheets.Item("de-de")

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/41709664619ba1ce34380b90cac491f7da408def/e2e"

$file1Md   = "df00d587-3ebe-43c4-8d41-16042e68dac1.md"
$file2Md   = "f7f91132-247a-4179-abb0-d9a421852111.md"
$file1Url  = "$repoBase/$file1Md"
$file2Url  = "$repoBase/$file2Md"

$statusHandedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet: update status text for both locales/rows
# ---------------------------------------------------------------------------
$ovr.Range("E2").Value = $statusHandedBack
$ovr.Range("F2").Value = $statusHandedBack
$ovr.Range("E3").Value = $statusHandedBack
$ovr.Range("F3").Value = $statusHandedBack

# ---------------------------------------------------------------------------
# zh-cn sheet: Latest Target File (I) / Latest Handback File (J) / Latest
# Handback DateTime (K)
# ---------------------------------------------------------------------------
$zhcn.Range("I2").Value = $file1Md
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $file1Url, "", "", $file1Md) | Out-Null
$zhcn.Range("J2").Value = "df00d587-3ebe-43c4-8d41-16042e68dac1.f6da3be473518da4571bb2ee44ff6fee3244cf68.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-26 04:24:13"

$zhcn.Range("I3").Value = $file2Md
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $file2Url, "", "", $file2Md) | Out-Null
$zhcn.Range("J3").Value = "f7f91132-247a-4179-abb0-d9a421852111.5003e3ff783e31ef79fed2c4f00fdcbf9fa2c3d6.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-26 04:24:13"

# ---------------------------------------------------------------------------
# de-de sheet: Latest Target File (I) / Latest Handback File (J) / Latest
# Handback DateTime (K)
# ---------------------------------------------------------------------------
$dede.Range("I2").Value = $file1Md
$dede.Hyperlinks.Add($dede.Range("I2"), $file1Url, "", "", $file1Md) | Out-Null
$dede.Range("J2").Value = "df00d587-3ebe-43c4-8d41-16042e68dac1.f6da3be473518da4571bb2ee44ff6fee3244cf68.de-de.xlf"
$dede.Range("K2").Value = "2016-08-26 04:24:20"

$dede.Range("I3").Value = $file2Md
$dede.Hyperlinks.Add($dede.Range("I3"), $file2Url, "", "", $file2Md) | Out-Null
$dede.Range("J3").Value = "f7f91132-247a-4179-abb0-d9a421852111.5003e3ff783e31ef79fed2c4f00fdcbf9fa2c3d6.de-de.xlf"
$dede.Range("K3").Value = "2016-08-26 04:24:20"

# ---------------------------------------------------------------------------
# Column widths: widen columns now holding the longer handback values
# ---------------------------------------------------------------------------
$ovr.Columns.Item(5).ColumnWidth  = 29.14   # E (zh-cn status)
$ovr.Columns.Item(6).ColumnWidth  = 29.14   # F (de-de status)

$zhcn.Columns.Item(3).ColumnWidth  = 29.14  # C (Status)
$zhcn.Columns.Item(9).ColumnWidth  = 39.14  # I (Latest Target File)
$zhcn.Columns.Item(10).ColumnWidth = 39.14  # J (Latest Handback File)

$dede.Columns.Item(3).ColumnWidth  = 29.14  # C (Status)
$dede.Columns.Item(9).ColumnWidth  = 39.14  # I (Latest Target File)
$dede.Columns.Item(10).ColumnWidth = 39.14  # J (Latest Handback File)
